$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.2

# Row 10
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57

# Row 11
$ws.Range("G11").Value = 2.15
$ws.Range("I11").Value = 3.7
$ws.Range("J11").Value = 3
$ws.Range("L11").Value = 4.33
$ws.Range("X11").Value = 9.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 19
$ws.Range("AF11").Value = 67
$ws.Range("AG11").Value = 9
$ws.Range("AH11").Value = 17
$ws.Range("AJ11").Value = 41
$ws.Range("AN11").Value = 4
$ws.Range("AR11").Value = 67
$ws.Range("AW11").Value = 5.5
$ws.Range("AZ11").Value = 81
